# edit.ps1 - applies the "How Psalms commentary works" edits via Word COM-interop.
$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Merge the "For Journal Editors and Academic Reviewers" paragraph into
#    the "Date: 2025-10-19" paragraph (removing the first paragraph's own
#    mark), and bump the date from 2025-10-19 to 2025-10-21, keeping the
#    date text split across two bold runs ("Date: 2025-10-" + "21").
# -----------------------------------------------------------------------
$editorsPara = $null
$dateParaIdx = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "For Journal Editors and Academic Reviewers*") {
        $editorsPara = $p
    }
    if ($t -like "Date: 2025-10-19*") {
        $dateParaIdx = $i
    }
    $i++
}

$nextPara = $d.Paragraphs($dateParaIdx)
$delRange = $d.Range($editorsPara.Range.Start, $nextPara.Range.Start)
$delRange.Delete()

# Locate the merged "Date: ..." paragraph again (index shifted by -1).
$dateP = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Date: 2025-10-19*") {
        $dateP = $p
    }
}

$dateStart = $dateP.Range.Start
$visibleLen = $dateP.Range.Text.Length - 1   # exclude trailing paragraph mark
$contentRange = $d.Range($dateStart, $dateStart + $visibleLen)
$contentRange.Delete()

$insPoint = $d.Range($dateStart, $dateStart)
$insPoint.InsertAfter("Date: 2025-10-21")

# Split into two runs ("Date: 2025-10-" / "21") while re-asserting Bold on
# each sub-range so the save path keeps them as distinct <w:r> elements.
$run1 = $d.Range($dateStart, $dateStart + 14)
$run1.Bold = 1
$run2 = $d.Range($dateStart + 14, $dateStart + 16)
$run2.Bold = 1

# -----------------------------------------------------------------------
# 2) Prefix the "Establishes the big picture..." paragraph (Stage 1 body)
#    with a new leading run: "A chapter of Psalms is fed to "
# -----------------------------------------------------------------------
$stage1Body = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Establishes the big picture*") {
        $stage1Body = $p
    }
}
$insertPos = $stage1Body.Range.Start
$insertTarget = $d.Range($insertPos, $insertPos)
$insertTarget.InsertBefore("A chapter of Psalms is fed to ")

$prefixLen = "A chapter of Psalms is fed to ".Length
$prefixRange = $d.Range($insertPos, $insertPos + $prefixLen)
$prefixRange.Bold = 1
$prefixRange.Bold = 0

# -----------------------------------------------------------------------
# 3) Collapse the runs of four paragraphs (which alternate English/Hebrew
#    text across many small <w:r> elements) down into a single run each,
#    without changing the visible text.
# -----------------------------------------------------------------------
function Flatten-Paragraph([string]$marker) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "$marker*") {
            $target = $p
        }
    }
    $start = $target.Range.Start
    $len = $target.Range.Text.Length - 1
    $text = $target.Range.Text.Substring(0, $len)

    $full = $d.Range($start, $start + $len)
    $full.Delete()

    $ins = $d.Range($start, $start)
    $ins.InsertAfter($text)
}

Flatten-Paragraph "The system’s analysis of Psalm 1 demonstrates"
Flatten-Paragraph "For instance, in verse 1, the sequence"
Flatten-Paragraph "The system’s analysis of the tree metaphor in Psalm 1:3"
Flatten-Paragraph "The system’s database reveals that"

Write-Output "done"
